# Bugfixed evaluation and simulated rt_data for components
# Rewrites the YoY forecast vector table (A2:E53) with refreshed simulated
# values and one additional leading observation (row 2 / FY2007).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# date_of_forecast (A), y_0 (B), y_0_forecast (C), y_1 (D), y_1_forecast (E)
$rows = @(
    ,@(2, 39400, 2007, 1.75539628881467, 2008, $null)
    ,@(3, 39583, 2008, $null, 2009, $null)
    ,@(4, 39765, 2008, 2.213911448916162, 2009, $null)
    ,@(5, 39948, 2009, $null, 2010, $null)
    ,@(6, 40130, 2009, 2.533533936850563, 2010, $null)
    ,@(7, 40310, 2010, 2.208165160720954, 2011, 1.903751357432193)
    ,@(8, 40494, 2010, 2.088987486264915, 2011, 1.485473821631844)
    ,@(9, 40676, 2011, 1.614140618728332, 2012, 1.770808585446004)
    ,@(10, 40862, 2011, 1.212544822741002, 2012, 1.799394172339341)
    ,@(11, 41044, 2012, 1.625793900975747, 2013, 1.586821460965226)
    ,@(12, 41228, 2012, 1.196776590518644, 2013, 1.2151583353186)
    ,@(13, 41409, 2013, 0.8049364973309325, 2014, 1.421244400332)
    ,@(14, 41592, 2013, 0.4712609263772594, 2014, 1.107727073902187)
    ,@(15, 41774, 2014, 0.5775251578155283, 2015, 1.341244385861273)
    ,@(16, 41957, 2014, 0.8783377572271434, 2015, 1.612081704302182)
    ,@(17, 42137, 2015, 1.901826580533572, 2016, 1.53605963063923)
    ,@(18, 42321, 2015, 2.29066283401107, 2016, 2.221748592150097)
    ,@(19, 42503, 2016, 2.590339257583607, 2017, 1.672072534917302)
    ,@(20, 42689, 2016, 4.109890522944348, 2017, 2.932944072183674)
    ,@(21, 42867, 2017, 1.713587272940131, 2018, 1.721854626734953)
    ,@(22, 43053, 2017, 1.336316831462692, 2018, 1.104283769064729)
    ,@(23, 43145, 2018, 1.808022822788802, 2019, 1.867774135387434)
    ,@(24, 43235, 2018, 1.05432456490544, 2019, 1.415552619392124)
    ,@(25, 43326, 2018, 1.299469465444592, 2019, 1.645976944955962)
    ,@(26, 43418, 2018, 1.197912858979611, 2019, 1.649865498505276)
    ,@(27, 43510, 2019, 2.247656020455691, 2020, 1.934175841213626)
    ,@(28, 43600, 2019, 1.566023898188384, 2020, 1.644188696416427)
    ,@(29, 43691, 2019, 1.75655962297816, 2020, 2.050351917667315)
    ,@(30, 43783, 2019, 1.727537197898665, 2020, 2.284828905445169)
    ,@(31, 43875, 2020, 2.684967757027334, 2021, 2.075491449101596)
    ,@(32, 43966, 2020, 2.155932165770968, 2021, 1.805141163113122)
    ,@(33, 44068, 2020, 2.980209378995857, 2021, 2.653391228709334)
    ,@(34, 44159, 2020, 3.647228437274408, 2021, 3.474365686630398)
    ,@(35, 44251, 2021, 1.906805170974435, 2022, 1.938263709207333)
    ,@(36, 44341, 2021, 2.443967114785739, 2022, 2.026008136667135)
    ,@(37, 44432, 2021, 1.954146674711188, 2022, 1.667670056759474)
    ,@(38, 44525, 2021, 2.777797690741424, 2022, 1.742844348069261)
    ,@(39, 44617, 2022, 1.467237762893392, 2023, 1.862063279188941)
    ,@(40, 44706, 2022, 0.388123216496683, 2023, 1.819907598678561)
    ,@(41, 44798, 2022, 2.69102598245059, 2023, 3.239034933968399)
    ,@(42, 44890, 2022, 0.6994919452575576, 2023, 0.5651273241891186)
    ,@(43, 44981, 2023, 0.2809429127725194, 2024, 1.721404396148163)
    ,@(44, 45071, 2023, -2.811030211656218, 2024, 0.8407670860975047)
    ,@(45, 45163, 2023, -1.669605379075589, 2024, 0.6305126186323617)
    ,@(46, 45254, 2023, -1.432689847121871, 2024, 0.4518870186319468)
    ,@(47, 45345, 2024, 1.069839250900739, 2025, 1.634674340565567)
    ,@(48, 45436, 2024, 1.250641979737566, 2025, 1.466559393695466)
    ,@(49, 45534, 2024, 1.780300968358017, 2025, 1.762346671645298)
    ,@(50, 45618, 2024, 2.033479419175133, 2025, 1.959987726090251)
    ,@(51, 45713, 2025, 2.97447584856072, 2026, 1.953339169714385)
    ,@(52, 45800, 2025, 2.302179720973463, 2026, 1.805984941845473)
    ,@(53, 45891, 2025, 2.481068287768839, 2026, 1.908500198348873)
)

foreach ($row in $rows) {
    $r  = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]   # A: date_of_forecast
    $ws.Cells.Item($r, 2).Value = $row[2]   # B: y_0
    $ws.Cells.Item($r, 3).Value = $row[3]   # C: y_0_forecast
    $ws.Cells.Item($r, 4).Value = $row[4]   # D: y_1
    $ws.Cells.Item($r, 5).Value = $row[5]   # E: y_1_forecast
}

# Row 53 is brand new: give its date cell (A53) the same date/time number
# format used by the rest of column A (copy format only, values already set above).
$ws.Cells.Item(52, 1).Copy() | Out-Null
$ws.Cells.Item(53, 1).PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

Write-Host "done"
